$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking price strings (e.g. "1.004", "27.374.08")
# that must stay as literal text exactly as scraped (same as original inlineStr cells).
# Force text format before assignment, then strip the style back off so the
# cell keeps no explicit style index (matching the source file).

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.374.08'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.68%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.881.35'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.12%  '
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '307.95'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.61%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.004'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.14%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5201'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3772'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.99%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07171'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.65%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.83'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.15%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.8893'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.18%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.906.12'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.08%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07597'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.71%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.359'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.21%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '89.95'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.53%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.007'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.33%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008590'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.20%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '14.20'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.19%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.003'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.10%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '27.419.99'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.63%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.087'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.30%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.128.44'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.72%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.65'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.61%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.509'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.45%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '152.45'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.23%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.850'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.33%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.09'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.53%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.159'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.11%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '113.07'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.93%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.764'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.16%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.709'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.57%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09056'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.69%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05189'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.29%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.110'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.23%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.189'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.20%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7571'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.71%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02051'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.18%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.521'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.78%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.058'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.45%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.087'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.07%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5487'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.02%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.657'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.39%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '115.81'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.96%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.487'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.48%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.1485'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.56%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4707'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.86%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.23'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.84%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.003'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.14%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.578'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.56%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '65.54'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.76%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '36.61'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.19%  '
